$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: change looked-up fruit to ŞEFTALİ, add COUNTIF formula
$ws.Range("G9").Value = "ŞEFTALİ"
$ws.Range("H9").Formula = "=COUNTIF(C4:E20,G9)"

# Row 14: change looked-up fruit to ARMUT, add SUMIF formulas
$ws.Range("G14").Value = "ARMUT"
$ws.Range("H14").Formula = "=SUMIF(C4:E20,G14,D4:D20)"
$ws.Range("I14").Formula = "=SUMIF(C4:E20,G14,E4:E20)"

# Row 17: Numara
$ws.Range("H17").Value = 20215070019

# Row 18: Ad Soyad
$ws.Range("H18").Value = "KÜBRA ÇABUK"

# Row 19: Bölüm
$ws.Range("H19").Value = "YBS"

# Update view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("H19:J19").Select()
